# Testing of valid user registration
# Adds two new worksheets ("registration" and "Sheet1") with test-data for
# the registration flow, re-points the "newsletter" sheet's selection, and
# makes "registration" the active tab.

$wb = $excel.ActiveWorkbook

# Right single quotation mark (U+2019) used in "Women's fashion" /
# "Men's fashion" - build it explicitly so it matches the existing
# shared-string text byte-for-byte regardless of script file encoding.
$rsquo = [char]0x2019
$womensFashion = "Women" + $rsquo + "s fashion"
$mensFashion = "Men" + $rsquo + "s fashion"

# ----------------------------------------------------------------------
# 1. "newsletter" sheet: move the selection from A5 to B2:B5 (no longer
#    the active/selected tab once the new sheets exist).
# ----------------------------------------------------------------------
$newsletter = $wb.Worksheets.Item("newsletter")
$newsletter.Range("B2:B5").Select()

# ----------------------------------------------------------------------
# 2. Add "registration" sheet right after "newsletter".
# ----------------------------------------------------------------------
$registration = $wb.Worksheets.Add($null, $newsletter)
$registration.Name = "registration"

$registration.Range("A1").Value = "First Name"
$registration.Range("B1").Value = "Last Name"
$registration.Range("C1").Value = "Email Address"
$registration.Range("D1").Value = "Password"
$registration.Range("E1").Value = "Category"
$registration.Range("A1:E1").Font.Size = 14
$registration.Range("A1:E1").Interior.Color = 65535
$registration.Rows.Item(1).RowHeight = 18

$registration.Range("A2").Value = "Washi"
$registration.Range("B2").Value = "Sundar"
$registration.Range("C2").Value = "washingtonsundar1800@gmail.com"
$registration.Hyperlinks.Add($registration.Range("C2"), "mailto:washingtonsundar1800@gmail.com")
$registration.Range("C2").Style = "Hyperlink"
$registration.Range("D2").Value = "zalandopwdpuma"
$registration.Range("E2").Value = $mensFashion
$registration.Range("E2").Font.Size = 14
$registration.Range("E2").Font.Color = 2236962
$registration.Rows.Item(2).RowHeight = 18

# Column widths (registration sheet) - closest values reachable through
# the interactive "set column width" pixel grid.
$registration.Columns.Item(1).ColumnWidth = 13.34
$registration.Columns.Item(2).ColumnWidth = 19.34
$registration.Columns.Item(3).ColumnWidth = 13.84
$registration.Columns.Item(4).ColumnWidth = 32.34
$registration.Columns.Item(5).ColumnWidth = 15.84
$registration.Columns.Item(6).ColumnWidth = 20.18

$registration.PageSetup.Orientation = 1

# ----------------------------------------------------------------------
# 3. Add "Sheet1" sheet right after "registration".
# ----------------------------------------------------------------------
$sheet1new = $wb.Worksheets.Add($null, $registration)
$sheet1new.Name = "Sheet1"

$sheet1new.Range("A1").Value = "P"
$sheet1new.Range("B1").Value = "hikip"
$sheet1new.Range("C1").Value = "raj"
$sheet1new.Range("D1").Value = "hikip75728@gmail.com"
$sheet1new.Hyperlinks.Add($sheet1new.Range("D1"), "mailto:hikip75728@gmail.com")
$sheet1new.Range("D1").Style = "Hyperlink"
$sheet1new.Range("E1").Value = "zalandopwdpuma"
$sheet1new.Range("F1").Value = $womensFashion
$sheet1new.Range("F1").Font.Size = 14
$sheet1new.Range("F1").Font.Color = 2236962
$sheet1new.Rows.Item(1).RowHeight = 18

$sheet1new.Range("A2").Value = "N"
$sheet1new.Range("C2").Value = "qwerty"
$sheet1new.Hyperlinks.Add($sheet1new.Range("D2"), "mailto:fashionzalandotest123@gmail.com", "", "", "fashionzalandotest123@gmail.com")
$sheet1new.Range("D2").Value = 123
$sheet1new.Range("D2").Style = "Hyperlink"
$sheet1new.Range("E2").Value = "zalandopwdpuma"
$sheet1new.Range("F2").Value = $mensFashion
$sheet1new.Range("F2").Font.Size = 14
$sheet1new.Range("F2").Font.Color = 2236962
$sheet1new.Rows.Item(2).RowHeight = 18

$sheet1new.Range("A3").Value = "P"
$sheet1new.Range("B3").Value = "testmail"
$sheet1new.Range("C3").Value = "zalando"
$sheet1new.Range("D3").Value = "testmailzalando123@gmail.com"
$sheet1new.Hyperlinks.Add($sheet1new.Range("D3"), "mailto:testmailzalando123@gmail.com")
$sheet1new.Range("D3").Style = "Hyperlink"
$sheet1new.Range("E3").Value = "zalandopwdpuma"
$sheet1new.Range("F3").Value = $womensFashion
$sheet1new.Range("F3").Font.Size = 14
$sheet1new.Range("F3").Font.Color = 2236962
$sheet1new.Rows.Item(3).RowHeight = 18

# ----------------------------------------------------------------------
# 4. Selections / active tab: "Sheet1" selects its full data range first,
#    then "registration" is activated and selects its full data range so
#    it ends up the active sheet (matches activeTab=1 / tabSelected).
# ----------------------------------------------------------------------
$sheet1new.Range("A1:F3").Select()
$registration.Activate()
$registration.Range("A1:E2").Select()

Write-Host "Registration test-data sheets added."
